$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell - reuse the same style/format as the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column
$values = @(0, 1, 0, 1, 1, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
